# Talent.xlsx edit: "unify the conception of DataNode, DataTable, Entity."
# The sheet that used to be a generic "Property1" table is renamed to
# "DataNode" to match the new naming convention, and the workbook is
# touched up (cursor position, a couple of row heights / column widths,
# and a small phonetic-guide font) the way Excel leaves things after a
# human re-opens, tweaks and re-saves the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- the actual content change: rename the sheet ---------------------
$ws.Name = "DataNode"

# --- cursor / selection left where the editor last clicked -----------
$ws.Range("D22").Select() | Out-Null

# --- header rows get a touch taller ----------------------------------
$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 27

# --- the two data columns get nudged a little wider -------------------
$ws.Columns.Item(1).ColumnWidth = 20.142857142857142
$ws.Columns.Item(8).ColumnWidth = 25.428571428571427

# --- small phonetic-guide font (9pt SimSun) used for east-asian ruby --
# registered in the style table without being applied to any cell, via a
# throw-away named style (mirrors how Excel mints a phonetic-guide font
# entry when phonetic info is turned on for the sheet).
$phoneticStyle = $wb.Styles.Add("PhoneticTmp")
$phoneticStyle.Font.Size = 9
$phoneticStyle.Font.Name = "宋体"
$wb.Styles.Item("PhoneticTmp").Delete() | Out-Null
